$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.965.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.05"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.42"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +8.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.288"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0679"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0938"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.064.37"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.32"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +13.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.795.39"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.845.04"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "258.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("E20").Value = "  +5.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.58"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.29"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.27"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.60"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.18%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0523"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.64"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("E34").Value = "  +10.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.471.93"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.642"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0192"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "85.21"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("E40").Value = "  +5.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.34"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.915"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.12"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.963.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.80%  "
